$wb = $excel.ActiveWorkbook

# --- Sheet "Fakturácia" (sheet1) ---
$ws1 = $wb.Worksheets.Item("Fakturácia")

# Clear the leftover "next period" rows: invoice-number values in A28:A31
# and the stray subtotal formulas in I28:K28 (keeps the cell styles).
[void]$ws1.Range("A28:A31").ClearContents()
[void]$ws1.Range("I28:K28").ClearContents()

# --- Sheet "Výnos celkom" (sheet2) ---
$ws2 = $wb.Worksheets.Item("Výnos celkom")
[void]$ws2.Range("E31").Select()

# Make "Fakturácia" the active sheet/tab again, with rows 28:33 selected
# (topLeftCell A16, selection A28:XFD33 == whole rows 28-33).
$ws1.Activate()
[void]$ws1.Range("A28:XFD33").Select()
